$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at K:L. This shifts the existing "Notes" column
# (old K) to M, and leaves the existing I ("delta.TPC") / J ("delta.model")
# columns untouched in place.
$ws.Columns("K:L").Insert()

# The existing I/J columns currently hold the *proportional* differences
# ((F-E)/E and (H-G)/G). Those formulas are being relocated to the new K/L
# columns (with new header labels below), while I/J are repurposed to hold
# the plain (absolute) differences F-E and H-G.
for ($r = 2; $r -le 30; $r++) {
    $iFormula = $ws.Range("I$r").Formula()
    $jFormula = $ws.Range("J$r").Formula()

    if ($iFormula -ne $null -and $iFormula -ne "") {
        $ws.Range("K$r").Formula = $iFormula
    }
    if ($jFormula -ne $null -and $jFormula -ne "") {
        $ws.Range("L$r").Formula = $jFormula
    }

    $ws.Range("I$r").Formula = "=F$r-E$r"
    $ws.Range("J$r").Formula = "=H$r-G$r"
}

# New header labels for the relocated proportional-difference columns.
$ws.Range("K1").Value = "delta.prop.TPC"
$ws.Range("L1").Value = "delta.prop.model"

# Give every formula cell in I:L a consistent numeric display format.
$ws.Range("I2:L30").NumberFormat = "0.000"

# Restore the selection to reflect where the edits were focused.
$ws.Range("I6").Select()
